$d = $word.ActiveDocument

# The document has an image paragraph, followed by 20 empty (formatting-only)
# paragraphs, followed by the final bold "Constantia" paragraph. The edit
# removes 17 of those 20 empty trailing paragraphs, leaving just 3 of them
# before the final paragraph (a manual clean-up pass before final submission).

$startPara = 51
$endPara = 67

$rStart = $d.Paragraphs.Item($startPara).Range.Start
$rEnd = $d.Paragraphs.Item($endPara).Range.End

$rng = $d.Range($rStart, $rEnd)
$rng.Delete()
